$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.257.44'
$ws.Range("E2").Value = '  +2.44%  '
$ws.Range("D3").Value = '1.877.63'
$ws.Range("E3").Value = '  +1.69%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.04'
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4308'
$ws.Range("E7").Value = '  +1.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3732'
$ws.Range("E8").Value = '  +2.73%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07409'
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8846'
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.13'
$ws.Range("E11").Value = '  +2.04%  '
$ws.Range("D12").Value = '1.950.98'
$ws.Range("E12").Value = '  +4.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.497'
$ws.Range("E13").Value = '  +3.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.636'
$ws.Range("E14").Value = '  +1.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06999'
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '81.33'
$ws.Range("E17").Value = '  +2.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009126'
$ws.Range("E18").Value = '  +2.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.63'
$ws.Range("E20").Value = '  +1.99%  '
$ws.Range("D21").Value = '28.369.71'
$ws.Range("E21").Value = '  +2.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.082'
$ws.Range("E22").Value = '  +2.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.95'
$ws.Range("E23").Value = '  +5.22%  '
$ws.Range("D24").Value = '2.229.86'
$ws.Range("E24").Value = '  +6.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.971'
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.56'
$ws.Range("E26").Value = '  +1.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.81'
$ws.Range("E27").Value = '  -0.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.392'
$ws.Range("E28").Value = '  +2.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.25'
$ws.Range("E29").Value = '  -3.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.866'
$ws.Range("E30").Value = '  -0.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08992'
$ws.Range("E31").Value = '  +1.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7893'
$ws.Range("E32").Value = '  +3.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.684'
$ws.Range("E33").Value = '  +2.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.179'
$ws.Range("E34").Value = '  +7.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.953'
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.127'
$ws.Range("E37").Value = '  +3.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05458'
$ws.Range("E38").Value = '  +2.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01966'
$ws.Range("E39").Value = '  +1.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.888'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5173'
$ws.Range("E41").Value = '  +1.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1686'
$ws.Range("E42").Value = '  +2.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.874'
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.667'
$ws.Range("E44").Value = '  +4.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.55'
$ws.Range("E45").Value = '  +2.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06607'
$ws.Range("E46").Value = '  +1.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4769'
$ws.Range("E47").Value = '  +0.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '106.20'
$ws.Range("E48").Value = '  +1.44%  '
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.656'
$ws.Range("E50").Value = '  +2.11%  '
$ws.Range("E51").Value = '  +5.60%  '
